$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert a new "Session 2.5 (SEMANTiCS)" subtitle paragraph (Heading5)
#    right before the existing "Time: Thursday, ..." (Heading4) paragraph,
#    which is the second paragraph of the document.
# ---------------------------------------------------------------------------
$timePara = $d.Paragraphs.Item(2)
$timePara.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs.Item(2)
$sessionXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="Heading5"/>
    <w:spacing w:after="80" w:before="240" w:lineRule="auto"/>
    <w:rPr/>
  </w:pPr>
  <w:bookmarkStart w:colFirst="0" w:colLast="0" w:name="_hezhfq48qg3k" w:id="1"/>
  <w:bookmarkEnd w:id="1"/>
  <w:r>
    <w:rPr>
      <w:rtl w:val="0"/>
    </w:rPr>
    <w:t xml:space="preserve">Session 2.5 (SEMANTiCS)</w:t>
  </w:r>
</w:p>
'@
$newPara.Range.InsertXML($sessionXml)

# ---------------------------------------------------------------------------
# 2. Replace the "Chair: TBA" placeholder with the actual chair information.
#    After the insertion above this paragraph moved from index 3 to index 4.
# ---------------------------------------------------------------------------
$chairPara = $d.Paragraphs.Item(4)
$chairPara.Range.Text = "Chair: Maribel Acosta, Assistant Professor, Ruhr University Bochum"

Write-Output "done"
